# Natmi following Dr Hou advice
# Update the LR-pairs sheet (Icam5 -> Itgal) to include the full set of
# sending/target cluster combinations (ECs, FAPs, sCs) instead of the
# single FAPs->FAPs row that was previously present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: ECs (sending) -> ECs (target) -------------------------------
$ws.Cells.Item(2, 1).Value  = "ECs"
$ws.Cells.Item(2, 2).Value  = "Icam5"
$ws.Cells.Item(2, 3).Value  = "Itgal"
$ws.Cells.Item(2, 4).Value  = "ECs"
$ws.Cells.Item(2, 5).Value  = 2
$ws.Cells.Item(2, 6).Value  = 0.6666666666666666
$ws.Cells.Item(2, 7).Value  = 2.064164333333333
$ws.Cells.Item(2, 8).Value  = 6.192493
$ws.Cells.Item(2, 9).Value  = 0.6667582237734649
$ws.Cells.Item(2, 10).Value = 0.6667582237734649
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 21.087087
$ws.Cells.Item(2, 14).Value = 63.261261
$ws.Cells.Item(2, 15).Value = 0.9808818221599021
$ws.Cells.Item(2, 16).Value = 0.9808818221599021
$ws.Cells.Item(2, 17).Value = 43.527212879297
$ws.Cells.Item(2, 18).Value = 391.744915913673
$ws.Cells.Item(2, 19).Value = 0.654011021475016
$ws.Cells.Item(2, 20).Value = 0.654011021475016

# ---- Row 3: ECs (sending) -> FAPs (target) -------------------------------
$ws.Cells.Item(3, 1).Value  = "ECs"
$ws.Cells.Item(3, 2).Value  = "Icam5"
$ws.Cells.Item(3, 3).Value  = "Itgal"
$ws.Cells.Item(3, 4).Value  = "FAPs"
$ws.Cells.Item(3, 5).Value  = 2
$ws.Cells.Item(3, 6).Value  = 0.6666666666666666
$ws.Cells.Item(3, 7).Value  = 2.064164333333333
$ws.Cells.Item(3, 8).Value  = 6.192493
$ws.Cells.Item(3, 9).Value  = 0.6667582237734649
$ws.Cells.Item(3, 10).Value = 0.6667582237734649
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3230143333333333
$ws.Cells.Item(3, 14).Value = 0.9690430000000001
$ws.Cells.Item(3, 15).Value = 0.01502525635066456
$ws.Cells.Item(3, 16).Value = 0.01502525635066456
$ws.Cells.Item(3, 17).Value = 0.6667546660221111
$ws.Cells.Item(3, 18).Value = 6.000791994199
$ws.Cells.Item(3, 19).Value = 0.01001821323611007
$ws.Cells.Item(3, 20).Value = 0.01001821323611007

# ---- Row 4: ECs (sending) -> sCs (target) --------------------------------
$ws.Cells.Item(4, 1).Value  = "ECs"
$ws.Cells.Item(4, 2).Value  = "Icam5"
$ws.Cells.Item(4, 3).Value  = "Itgal"
$ws.Cells.Item(4, 4).Value  = "sCs"
$ws.Cells.Item(4, 5).Value  = 2
$ws.Cells.Item(4, 6).Value  = 0.6666666666666666
$ws.Cells.Item(4, 7).Value  = 2.064164333333333
$ws.Cells.Item(4, 8).Value  = 6.192493
$ws.Cells.Item(4, 9).Value  = 0.6667582237734649
$ws.Cells.Item(4, 10).Value = 0.6667582237734649
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.08799
$ws.Cells.Item(4, 14).Value = 0.26397
$ws.Cells.Item(4, 15).Value = 0.00409292148943331
$ws.Cells.Item(4, 16).Value = 0.004092921489433309
$ws.Cells.Item(4, 17).Value = 0.18162581969
$ws.Cells.Item(4, 18).Value = 1.63463237721
$ws.Cells.Item(4, 19).Value = 0.002728989062338798
$ws.Cells.Item(4, 20).Value = 0.002728989062338798

# ---- Row 5: FAPs (sending) -> ECs (target) -------------------------------
$ws.Cells.Item(5, 1).Value  = "FAPs"
$ws.Cells.Item(5, 2).Value  = "Icam5"
$ws.Cells.Item(5, 3).Value  = "Itgal"
$ws.Cells.Item(5, 4).Value  = "ECs"
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 6).Value  = 1
$ws.Cells.Item(5, 7).Value  = 1.031657
$ws.Cells.Item(5, 8).Value  = 3.094971
$ws.Cells.Item(5, 9).Value  = 0.3332417762265351
$ws.Cells.Item(5, 10).Value = 0.3332417762265351
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 21.087087
$ws.Cells.Item(5, 14).Value = 63.261261
$ws.Cells.Item(5, 15).Value = 0.9808818221599021
$ws.Cells.Item(5, 16).Value = 0.9808818221599021
$ws.Cells.Item(5, 17).Value = 21.754640913159
$ws.Cells.Item(5, 18).Value = 195.791768218431
$ws.Cells.Item(5, 19).Value = 0.3268708006848861
$ws.Cells.Item(5, 20).Value = 0.3268708006848861

# ---- Row 6: FAPs (sending) -> FAPs (target) ------------------------------
$ws.Cells.Item(6, 1).Value  = "FAPs"
$ws.Cells.Item(6, 2).Value  = "Icam5"
$ws.Cells.Item(6, 3).Value  = "Itgal"
$ws.Cells.Item(6, 4).Value  = "FAPs"
$ws.Cells.Item(6, 5).Value  = 3
$ws.Cells.Item(6, 6).Value  = 1
$ws.Cells.Item(6, 7).Value  = 1.031657
$ws.Cells.Item(6, 8).Value  = 3.094971
$ws.Cells.Item(6, 9).Value  = 0.3332417762265351
$ws.Cells.Item(6, 10).Value = 0.3332417762265351
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3230143333333333
$ws.Cells.Item(6, 14).Value = 0.9690430000000001
$ws.Cells.Item(6, 15).Value = 0.01502525635066456
$ws.Cells.Item(6, 16).Value = 0.01502525635066456
$ws.Cells.Item(6, 17).Value = 0.3332399980836667
$ws.Cells.Item(6, 18).Value = 2.999159982753
$ws.Cells.Item(6, 19).Value = 0.005007043114554483
$ws.Cells.Item(6, 20).Value = 0.005007043114554483

# ---- Row 7: FAPs (sending) -> sCs (target) -------------------------------
$ws.Cells.Item(7, 1).Value  = "FAPs"
$ws.Cells.Item(7, 2).Value  = "Icam5"
$ws.Cells.Item(7, 3).Value  = "Itgal"
$ws.Cells.Item(7, 4).Value  = "sCs"
$ws.Cells.Item(7, 5).Value  = 3
$ws.Cells.Item(7, 6).Value  = 1
$ws.Cells.Item(7, 7).Value  = 1.031657
$ws.Cells.Item(7, 8).Value  = 3.094971
$ws.Cells.Item(7, 9).Value  = 0.3332417762265351
$ws.Cells.Item(7, 10).Value = 0.3332417762265351
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.08799
$ws.Cells.Item(7, 14).Value = 0.26397
$ws.Cells.Item(7, 15).Value = 0.00409292148943331
$ws.Cells.Item(7, 16).Value = 0.004092921489433309
$ws.Cells.Item(7, 17).Value = 0.09077549943
$ws.Cells.Item(7, 18).Value = 0.81697949487
$ws.Cells.Item(7, 19).Value = 0.001363932427094512
$ws.Cells.Item(7, 20).Value = 0.001363932427094511
